$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Prudential Financial, Inc. (PRU)
$ws.Range("E2").Value = 62.7
$ws.Range("G2").Value = 60
$ws.Range("K2").Value = 62.2
$ws.Range("M2").Value = "📈 매수 관찰 구간입니다."
$ws.Range("N2").Value = 85.87127175646313
$ws.Range("O2").Value = "🟢 완화적 (상승 우위)"

# Row 3 - UnitedHealth Group Incorporated (UNH)
$ws.Range("E3").Value = 53.8
$ws.Range("G3").Value = 50
$ws.Range("K3").Value = 59.2
$ws.Range("M3").Value = "⛔ 관망하십시오."
$ws.Range("N3").Value = 85.87127175646313
$ws.Range("O3").Value = "🟢 완화적 (상승 우위)"

# Row 4 - MetLife, Inc. (MET)
$ws.Range("E4").Value = 51.6
$ws.Range("G4").Value = 40
$ws.Range("K4").Value = 57.8
$ws.Range("M4").Value = "⛔ 관망하십시오."
$ws.Range("N4").Value = 85.87127175646313
$ws.Range("O4").Value = "🟢 완화적 (상승 우위)"

# Row 5 - American International Group, I (AIG)
$ws.Range("E5").Value = 45.8
$ws.Range("G5").Value = 30
$ws.Range("K5").Value = 54.8
$ws.Range("M5").Value = "⛔ 관망하십시오."
$ws.Range("N5").Value = 85.87127175646313
$ws.Range("O5").Value = "🟢 완화적 (상승 우위)"
